# Apply corrected partner data edits
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix mojibake text in A28 (double-encoded UTF-8 -> Latin-1 artifact)
$ws.Range("A28").Value = "Consejo Nacional de Investigaciones CientÃ­ficas y TÃ©cnicas"

# Rows 39, 72, 90: change all values in columns B:AK from 2 to 1
$rows = @(39, 72, 90)
foreach ($r in $rows) {
    $rng = $ws.Range("B$r`:AK$r")
    $rng.Value = 1
}
